$d = $word.ActiveDocument

$pairs = @(
    @("2024-04-01 Monday", "2024-04-02 Tuesday"),
    @("283×7=", "403×7="),
    @("698×8=", "898×2="),
    @("561×3=", "674×4="),
    @("169×9=", "280×4="),
    @("182×8=", "953×6="),
    @("258×6=", "777×4="),
    @("542×8=", "502×7="),
    @("178×7=", "309×2="),
    @("732×9=", "183×8="),
    @("126×3=", "624×7="),
    @("891×8=", "158×3="),
    @("646×2=", "594×4="),
    @("814×8=", "584×4="),
    @("153×6=", "212×4="),
    @("343×4=", "379×8="),
    @("985×2=", "188×6="),
    @("695×4=", "874×8="),
    @("365×7=", "368×2="),
    @("668×9=", "309×3="),
    @("541×2=", "783×8="),
    @("191×8=", "265×3="),
    @("901×5=", "825×6="),
    @("852×5=", "519×7="),
    @("946×6=", "330×3="),
    @("680×7=", "199×9=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
